$d = $word.ActiveDocument

# The two Pearson logo pictures (in footer1.xml / footer2.xml) were
# exported as "image2.png" but should be named "image1.png"; the two
# BTEC logo pictures (in header1.xml / header2.xml) were exported as
# "image1.jpg" but should be named "image2.jpg". Only the internal
# wp:docPr/@name and pic:cNvPr/@name labels change - the embedded
# relationship targets (media/image2.png, media/image1.jpg) stay the
# same, as do all other attributes (id, descr, etc).
#
# InlineShape does not expose a settable Name property via this object
# model, so rewrite the underlying part XML directly through the
# document-wide WordOpenXML round-trip (covers headers/footers too).

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml
